$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Sheet2")

# Insert a new row at row 5 (pushes existing rows 5-14 down to 6-15),
# inheriting formatting from the row above as Excel normally does.
$ws2.Rows.Item(5).Insert()

# Populate the new row: category label + the same SUMIFS formulas used
# by every other row in the table.
$ws2.Range("A5").Value = "支出:杂项"
$ws2.Range("B5").Formula = "=SUMIFS(Sheet1!`$G`$2:`$G`$65475,Sheet1!`$E`$2:`$E`$65475,`$A5,Sheet1!`$F`$2:`$F`$65475,B`$1)-SUMIFS(Sheet1!`$G`$2:`$G`$65475,Sheet1!`$D`$2:`$D`$65475,`$A5,Sheet1!`$F`$2:`$F`$65475,B`$1)"
$ws2.Range("C5").Formula = "=SUMIFS(Sheet1!`$G`$2:`$G`$65475,Sheet1!`$E`$2:`$E`$65475,`$A5,Sheet1!`$F`$2:`$F`$65475,C`$1)-SUMIFS(Sheet1!`$G`$2:`$G`$65475,Sheet1!`$D`$2:`$D`$65475,`$A5,Sheet1!`$F`$2:`$F`$65475,C`$1)"

# The conditional-formatting range that used to cover B2:C14 should now
# extend to the freshly-inserted row (B2:C15).
$fc = $ws2.Range("B2:C14").FormatConditions.Item(2)
$fc.ModifyAppliesToRange($ws2.Range("B2:C15"))

# Make Sheet2 the active tab and select B5, matching the saved UI state.
$ws2.Activate()
$ws2.Range("B5").Select()
